# Insert a new weekly price record for "Terminal Hortofrutícola Agro Chillán"
# (Apio / Americana (o) / Primera) at row 163 of the dataset, pushing the
# existing rows 163-171 down to 164-172 (dimension grows from R171 to R172).
#
# The new row duplicates the former row 163's values except for the
# reporting date (column D) and the reported volume (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 163..171 down to 164..172, leaving row 163 blank for the new record.
$ws.Rows.Item(163).Insert()

$ws.Range("A163").Value = 7
$ws.Range("B163").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C163").Value = "Ñuble"
$ws.Range("D163").Value = 44585
$ws.Range("E163").Value = 16
$ws.Range("F163").Value = 100112017
$ws.Range("G163").Value = "Apio"
$ws.Range("H163").Value = "Americana (o)"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 100
$ws.Range("K163").Value = 8000
$ws.Range("L163").Value = 8500
$ws.Range("M163").Value = 8250
$ws.Range("N163").Value = "`$/docena de matas"
$ws.Range("O163").Value = "Provincia del Elquí"
$ws.Range("P163").Value = 1375
$ws.Range("Q163").Value = 6
$ws.Range("R163").Value = "Hortaliza"
